$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E in the data range to be treated as text so that
# numeric-looking strings (e.g. "602.60", "68.967.28") are not auto-converted
# to numbers when assigned via .Value
$numRng = $ws.Range("D2:E51")
$numRng.NumberFormat = "@"

$ws.Range("D2").Value = "68.967.28"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "3.744.83"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "602.60"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").Value = "167.39"
$ws.Range("E6").Value = "  +1.14%  "
$ws.Range("D7").Value = "3.742.81"
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("E10").Value = "  +2.79%  "
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").Value = "37.93"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").Value = "0.0000249"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").Value = "4.379.44"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").Value = "3.752.45"
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("D17").Value = "68.993.92"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("E19").Value = "  -1.12%  "
$ws.Range("D20").Value = "17.25"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").Value = "10.86"
$ws.Range("E21").Value = "  +19.69%  "
$ws.Range("D22").Value = "492.22"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "0.726"
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("E24").Value = "  +8.62%  "
$ws.Range("D25").Value = "84.76"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").Value = "2.31"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("D27").Value = "12.35"
$ws.Range("E27").Value = "  +1.76%  "
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  +2.68%  "
$ws.Range("E31").Value = "  +5.06%  "
$ws.Range("E32").Value = "  +2.75%  "
$ws.Range("D33").Value = "31.55"
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("D34").Value = "3.891.04"
$ws.Range("E34").Value = "  +2.04%  "
$ws.Range("B35").Value = "RenzoRestakedETH"
$ws.Range("C35").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D35").Value = "3.681.53"
$ws.Range("E35").Value = "  +1.76%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.108"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("E38").Value = "  +2.46%  "
$ws.Range("D39").Value = "5.90"
$ws.Range("E39").Value = "  +3.49%  "
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("D41").Value = "0.325"
$ws.Range("E41").Value = "  +1.58%  "
$ws.Range("D42").Value = "2.97"
$ws.Range("E42").Value = "  +6.34%  "
$ws.Range("D43").Value = "430.42"
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").Value = "1.99"
$ws.Range("E45").Value = "  +2.10%  "
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "40.40"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").Value = "141.07"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "2.783.14"
$ws.Range("E50").Value = "  +2.36%  "
$ws.Range("E51").Value = "  +1.41%  "

# Restore the default "Normal" style on that range (removes the temporary text
# number format) so the resulting cell styling matches the original file
$numRng.Style = "Normal"
